$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Sun Jan 14 16:32:46 EST 2024"
$ws.Range("B3").Value = "Sun Jan 14 16:32:59 EST 2024"
$ws.Range("B5").Value = "Sun Jan 14 16:33:11 EST 2024"
